$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 0.026481028646230698
$ws.Range("C2").Value = 0.011453032493591309
$ws.Range("D2").Value = 0.007475130259990692
$ws.Range("F2").Value = 0.00015327485743910074
$ws.Range("J2").Value = 0.12750037014484406
$ws.Range("K2").Value = 1.4436825513839722
